$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2 through 10 from 45174 to 45175
for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45174) {
        $cell.Value = 45175
    }
}
